$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 30 (week 29) value: was placeholder 1, now 342
$ws.Range("B30").Value = 342

# Add new row 31 (week 30) with its case count
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 36
